$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 352, shifting existing rows 352:440 down to 354:442.
$ws.Rows("352:353").Insert()

# Fill the new row 352 with fresh data (Packham's Triumph, Primera).
$ws.Cells.Item(352, 1).Value = 4
$ws.Cells.Item(352, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(352, 3).Value = "Los Lagos"
$ws.Cells.Item(352, 4).Value = 44995
$ws.Cells.Item(352, 5).Value = 10
$ws.Cells.Item(352, 6).Value = "Fruta"
$ws.Cells.Item(352, 7).Value = 100104
$ws.Cells.Item(352, 8).Value = "Frutos de pepita"
$ws.Cells.Item(352, 9).Value = 100104005
$ws.Cells.Item(352, 10).Value = "Pera"
$ws.Cells.Item(352, 11).Value = "Packham's Triumph"
$ws.Cells.Item(352, 12).Value = "Primera"
$ws.Cells.Item(352, 13).Value = 400
$ws.Cells.Item(352, 14).Value = 16000
$ws.Cells.Item(352, 15).Value = 17000
$ws.Cells.Item(352, 16).Value = 16500
$ws.Cells.Item(352, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(352, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(352, 19).Value = 1100
$ws.Cells.Item(352, 20).Value = 15

# Fill the new row 353 with fresh data (Packham's Triumph, Segunda).
$ws.Cells.Item(353, 1).Value = 4
$ws.Cells.Item(353, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(353, 3).Value = "Los Lagos"
$ws.Cells.Item(353, 4).Value = 44995
$ws.Cells.Item(353, 5).Value = 10
$ws.Cells.Item(353, 6).Value = "Fruta"
$ws.Cells.Item(353, 7).Value = 100104
$ws.Cells.Item(353, 8).Value = "Frutos de pepita"
$ws.Cells.Item(353, 9).Value = 100104005
$ws.Cells.Item(353, 10).Value = "Pera"
$ws.Cells.Item(353, 11).Value = "Packham's Triumph"
$ws.Cells.Item(353, 12).Value = "Segunda"
$ws.Cells.Item(353, 13).Value = 200
$ws.Cells.Item(353, 14).Value = 14000
$ws.Cells.Item(353, 15).Value = 14000
$ws.Cells.Item(353, 16).Value = 14000
$ws.Cells.Item(353, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(353, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(353, 19).Value = 933
$ws.Cells.Item(353, 20).Value = 15
